$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row number, Date string (new), D, E, F, G, H
$rows = @(
    @{ R = 3;  A = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ R = 4;  A = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ R = 5;  A = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ R = 6;  A = "08-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 7;  A = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 8;  A = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 9;  A = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 10; A = "22-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 11; A = "25-08-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ R = 12; A = "29-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ R = 13; A = "01-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ R = 14; A = "05-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 15; A = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 16; A = "12-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 17; A = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 18; A = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 19; A = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 20; A = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ R = 21; A = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($row in $rows) {
    $r = $row.R
    $cellA = $ws.Cells.Item($r, 1)
    # Temporarily force text format so Excel doesn't auto-convert the
    # dd-mm-yyyy-looking string into a date serial, then restore the
    # cell's original style so no stray formatting is introduced.
    $origStyle = $cellA.Style
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.A
    $cellA.Style = $origStyle

    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
